$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,13

$arr[0,0] = 0.8521272959555688
$arr[0,1] = 0.04841290214839944
$arr[0,2] = 0.07758153008460056
$arr[0,3] = 0.06955353243393247
$arr[0,4] = 0
$arr[0,5] = 1.4727967649996
$arr[0,6] = 1.338123636665202
$arr[0,7] = 1.291858295027353
$arr[0,8] = 0
$arr[0,9] = 0.6317458505050411
$arr[0,10] = 0.2339524660021084
$arr[0,11] = 0
$arr[0,12] = 2.280027468995655

$arr[1,0] = 0.8082002511848998
$arr[1,1] = 0.043290728730355
$arr[1,2] = 0.07050867933646998
$arr[1,3] = 0.06913241108798118
$arr[1,4] = 0
$arr[1,5] = 1.465085432766955
$arr[1,6] = 1.339706944032599
$arr[1,7] = 1.29272454619926
$arr[1,8] = 0
$arr[1,9] = 0.5863677974402322
$arr[1,10] = 0.2264942255563085
$arr[1,11] = 0
$arr[1,12] = 2.297769912496562

$arr[2,0] = 0.7816718013819184
$arr[2,1] = 0.0401251948017034
$arr[2,2] = 0.06620255808185505
$arr[2,3] = 0.06890477234951398
$arr[2,4] = 0
$arr[2,5] = 1.461144629844185
$arr[2,6] = 1.34123595468759
$arr[2,7] = 1.293836636706871
$arr[2,8] = 0
$arr[2,9] = 0.5588119837284751
$arr[2,10] = 0.2220486068125638
$arr[2,11] = 0
$arr[2,12] = 2.309342281414736

$arr[3,0] = 0.7709728964556462
$arr[3,1] = 0.03883000810346005
$arr[3,2] = 0.06445695025364273
$arr[3,3] = 0.06881979723076981
$arr[3,4] = 0
$arr[3,5] = 1.459738156722054
$arr[3,6] = 1.341999086471134
$arr[3,7] = 1.294435675269952
$arr[3,8] = 0
$arr[3,9] = 0.5476599543034979
$arr[3,10] = 0.2202706265379106
$arr[3,11] = 0
$arr[3,12] = 2.31422861629693

$arr[4,0] = 0.7692031032493674
$arr[4,1] = 0.03861462741920718
$arr[4,2] = 0.06416764614533577
$arr[4,3] = 0.06880615798644563
$arr[4,4] = 0
$arr[4,5] = 1.459516652066924
$arr[4,6] = 1.342134263150442
$arr[4,7] = 1.294543953023741
$arr[4,8] = 0
$arr[4,9] = 0.5458128338917732
$arr[4,10] = 0.2199774266201757
$arr[4,11] = 0
$arr[4,12] = 2.315050283483281

$arr[5,0] = 0.7815270597616006
$arr[5,1] = 0.04010774860563515
$arr[5,2] = 0.06617897912306603
$arr[5,3] = 0.06890359479031538
$arr[5,4] = 0
$arr[5,5] = 1.461124854439291
$arr[5,6] = 1.341245679477609
$arr[5,7] = 1.293844125062726
$arr[5,8] = 0
$arr[5,9] = 0.5586612708677308
$arr[5,10] = 0.2220244921084031
$arr[5,11] = 0
$arr[5,12] = 2.309407490022537

$arr[6,0] = 0.8368894683886481
$arr[6,1] = 0.0466509989896764
$arr[6,2] = 0.07513517989465868
$arr[6,3] = 0.06940191658232919
$arr[6,4] = 0
$arr[6,5] = 1.469972967684029
$arr[6,6] = 1.338553996375836
$arr[6,7] = 1.292036525452289
$arr[6,8] = 0
$arr[6,9] = 0.6160359396183992
$arr[6,10] = 0.2313531020340065
$arr[6,11] = 0
$arr[6,12] = 2.286004195974591

$arr[7,0] = 0.9489640445728185
$arr[7,1] = 0.05932251597779725
$arr[7,2] = 0.09299177029112116
$arr[7,3] = 0.07062417327765402
$arr[7,4] = 0
$arr[7,5] = 1.493636740984982
$arr[7,6] = 1.337694300101987
$arr[7,7] = 1.29309896681081
$arr[7,8] = 0
$arr[7,9] = 0.7309807267947122
$arr[7,10] = 0.2507088375137414
$arr[7,11] = 0
$arr[7,12] = 2.245497481070004

$arr[8,0] = 1.033447014211788
$arr[8,1] = 0.06853964864181705
$arr[8,2] = 0.1062951751219998
$arr[8,3] = 0.07167123660953578
$arr[8,4] = 0
$arr[8,5] = 1.514892652673836
$arr[8,6] = 1.33975873913343
$arr[8,7] = 1.296694950087925
$arr[8,8] = 0
$arr[8,9] = 0.81692548940606
$arr[8,10] = 0.2655803128647136
$arr[8,11] = 0
$arr[8,12] = 2.219026195927441

$arr[9,0] = 1.072346778152848
$arr[9,1] = 0.07271372020235845
$arr[9,2] = 0.1123884335528942
$arr[9,3] = 0.07217988874623771
$arr[9,4] = 0
$arr[9,5] = 1.525407873127563
$arr[9,6] = 1.34128399217866
$arr[9,7] = 1.298943820878428
$arr[9,8] = 0
$arr[9,9] = 0.8563516518578353
$arr[9,10] = 0.2724878531262505
$arr[9,11] = 0
$arr[9,12] = 2.207698999510569

$arr[10,0] = 1.087144294766233
$arr[10,1] = 0.07429169016843673
$arr[10,2] = 0.1147018221444966
$arr[10,3] = 0.07237714401248496
$arr[10,4] = 0
$arr[10,5] = 1.52951166710622
$arr[10,6] = 1.341945877463445
$arr[10,7] = 1.299883665557253
$arr[10,8] = 0
$arr[10,9] = 0.8713287309225848
$arr[10,10] = 0.2751240660321059
$arr[10,11] = 0
$arr[10,12] = 2.203512587990097

$arr[11,0] = 1.083954408536897
$arr[11,1] = 0.07395196362514866
$arr[11,2] = 0.1142033248580674
$arr[11,3] = 0.07233445537260508
$arr[11,4] = 0
$arr[11,5] = 1.528622415209213
$arr[11,6] = 1.341799578787146
$arr[11,7] = 1.299677326938493
$arr[11,8] = 0
$arr[11,9] = 0.8681010496769375
$arr[11,10] = 0.2745553997431927
$arr[11,11] = 0
$arr[11,12] = 2.204409624964086

$arr[12,0] = 1.073562837176496
$arr[12,1] = 0.0728435937707701
$arr[12,2] = 0.1125786368690598
$arr[12,3] = 0.07219602412767046
$arr[12,4] = 0
$arr[12,5] = 1.525743049875388
$arr[12,6] = 1.341336756157347
$arr[12,7] = 1.299019373350376
$arr[12,8] = 0
$arr[12,9] = 0.8575828789965954
$arr[12,10] = 0.2727043255753472
$arr[12,11] = 0
$arr[12,12] = 2.20735251658693

$arr[13,0] = 1.067206419326908
$arr[13,1] = 0.07216434026312868
$arr[13,2] = 0.1115842516334737
$arr[13,3] = 0.07211183491706308
$arr[13,4] = 0
$arr[13,5] = 1.52399523996263
$arr[13,6] = 1.341064243578472
$arr[13,7] = 1.298627852782161
$arr[13,8] = 0
$arr[13,9] = 0.8511463449161454
$arr[13,10] = 0.2715731558323711
$arr[13,11] = 0
$arr[13,12] = 2.209168535331663

$arr[14,0] = 1.030914217016687
$arr[14,1] = 0.06826648824555548
$arr[14,2] = 0.1058978043243428
$arr[14,3] = 0.07163864482515336
$arr[14,4] = 0
$arr[14,5] = 1.514222502095862
$arr[14,6] = 1.339670856375449
$arr[14,7] = 1.296560324821996
$arr[14,8] = 0
$arr[14,9] = 0.8143555154678381
$arr[14,10] = 0.2651317545875287
$arr[14,11] = 0
$arr[14,12] = 2.219780851841762

$arr[15,0] = 1.008769799873903
$arr[15,1] = 0.06587049033086601
$arr[15,2] = 0.1024200038496161
$arr[15,3] = 0.07135663309519913
$arr[15,4] = 0
$arr[15,5] = 1.50844407430597
$arr[15,6] = 1.338966194106177
$arr[15,7] = 1.295449042309286
$arr[15,8] = 0
$arr[15,9] = 0.7918698233063139
$arr[15,10] = 0.2612166406322558
$arr[15,11] = 0
$arr[15,12] = 2.226474345938755

$arr[16,0] = 0.9960769874765276
$arr[16,1] = 0.06449059861421347
$arr[16,2] = 0.1004235662353778
$arr[16,3] = 0.07119747180926339
$arr[16,4] = 0
$arr[16,5] = 1.505200069038125
$arr[16,6] = 1.338616058910503
$arr[16,7] = 1.294867558607528
$arr[16,8] = 0
$arr[16,9] = 0.7789676657024813
$arr[16,10] = 0.2589781771993813
$arr[16,11] = 0
$arr[16,12] = 2.230391529148683

$arr[17,0] = 0.991786997130248
$arr[17,1] = 0.06402308424902969
$arr[17,2] = 0.09974827557773835
$arr[17,3] = 0.0711441057292852
$arr[17,4] = 0
$arr[17,5] = 1.504115366280843
$arr[17,6] = 1.338506983711028
$arr[17,7] = 1.294680584907127
$arr[17,8] = 0
$arr[17,9] = 0.7746045481216015
$arr[17,10] = 0.2582225746645577
$arr[17,11] = 0
$arr[17,12] = 2.231729367191505

$arr[18,0] = 1.011122551793733
$arr[18,1] = 0.06612573195975813
$arr[18,2] = 0.1027898175954647
$arr[18,3] = 0.07138633871673861
$arr[18,4] = 0
$arr[18,5] = 1.509050958125016
$arr[18,6] = 1.339035496789421
$arr[18,7] = 1.295561368184103
$arr[18,8] = 0
$arr[18,9] = 0.7942602543374164
$arr[18,10] = 0.2616320236998888
$arr[18,11] = 0
$arr[18,12] = 2.225754849385901

$arr[19,0] = 1.076613276912951
$arr[19,1] = 0.07316922066611653
$arr[19,2] = 0.1130556838520249
$arr[19,3] = 0.07223655890113889
$arr[19,4] = 0
$arr[19,5] = 1.526585478447515
$arr[19,6] = 1.341470410194376
$arr[19,7] = 1.299210234553627
$arr[19,8] = 0
$arr[19,9] = 0.8606710394358856
$arr[19,10] = 0.2732474752122442
$arr[19,11] = 0
$arr[19,12] = 2.206485322477064

$arr[20,0] = 1.119805767426669
$arr[20,1] = 0.07775706448613562
$arr[20,2] = 0.1198000345678025
$arr[20,3] = 0.07281926613005751
$arr[20,4] = 0
$arr[20,5] = 1.53875600861312
$arr[20,6] = 1.343553172179213
$arr[20,7] = 1.302109389456582
$arr[20,8] = 0
$arr[20,9] = 0.9043496599726666
$arr[20,10] = 0.2809582000467827
$arr[20,11] = 0
$arr[20,12] = 2.194491739447379

$arr[21,0] = 1.096717482974441
$arr[21,1] = 0.07530984605170943
$arr[21,2] = 0.1161972306751125
$arr[21,3] = 0.07250579364336573
$arr[21,4] = 0
$arr[21,5] = 1.532195243612875
$arr[21,6] = 1.342396592036636
$arr[21,7] = 1.300514956611906
$arr[21,8] = 0
$arr[21,9] = 0.8810124071508483
$arr[21,10] = 0.2768319213101336
$arr[21,11] = 0
$arr[21,12] = 2.200837966192225

$arr[22,0] = 1.010058753304691
$arr[22,1] = 0.06601034470257616
$arr[22,2] = 0.1026226154740044
$arr[22,3] = 0.07137289955048232
$arr[22,4] = 0
$arr[22,5] = 1.508776342720239
$arr[22,6] = 1.339003993764209
$arr[22,7] = 1.295510406834694
$arr[22,8] = 0
$arr[22,9] = 0.7931794620710377
$arr[22,10] = 0.2614441904304385
$arr[22,11] = 0
$arr[22,12] = 2.226079918829555

$arr[23,0] = 0.918268978362903
$arr[23,1] = 0.05591114712187562
$arr[23,2] = 0.0881291257603749
$arr[23,3] = 0.07026731931833652
$arr[23,4] = 0
$arr[23,5] = 1.486557009185404
$arr[23,6] = 1.33745358925475
$arr[23,7] = 1.292317629538601
$arr[23,8] = 0
$arr[23,9] = 0.6996234449854342
$arr[23,10] = 0.2453585385739103
$arr[23,11] = 0
$arr[23,12] = 2.255878405946405

$ws.Range("B2:N25").Value = $arr
Write-Output "done"